# Apply the changes described by the diff to the FlashScore weekly-games
# workbook:
#   - small odds tweaks on row 2 (Q2, R2, U2, V2)
#   - AS3 updated on row 3
#   - the two Colombian "PRIMERA B" fixtures (old rows 4 & 5) are dropped
#   - the Paraguay fixture (old row 6) becomes the new row 4 (unchanged)
#   - the Spain fixture (old row 7) becomes the new row 5, with several of
#     its odds updated
#   - the sheet shrinks from A1:AS7 to A1:AS5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: a handful of odds changed -------------------------------------
$ws.Range("Q2").Value = 1.6
$ws.Range("R2").Value = 2.29
$ws.Range("U2").Value = 3.1
$ws.Range("V2").Value = 1.35

# --- Row 3: max payout (AS3) changed ---------------------------------------
$ws.Range("AS3").Value = 1000

# --- Drop the two old rows 4 & 5 (COLOMBIA - PRIMERA B fixtures). This
#     shifts the old row 6 (Paraguay) up to row 4 and the old row 7 (Spain)
#     up to row 5, matching the target layout exactly. ----------------------
$ws.Rows.Item(4).EntireRow.Delete() | Out-Null
$ws.Rows.Item(4).EntireRow.Delete() | Out-Null

# --- Row 5 (formerly row 7, SPAIN - LALIGA2) needs several odds refreshed --
$ws.Range("G5").Value = 2.2
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 3.2
$ws.Range("J5").Value = 3
$ws.Range("L5").Value = 4
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("S5").Value = 2.15
$ws.Range("T5").Value = 1.67
$ws.Range("W5").Value = 4
$ws.Range("X5").Value = 1.22
$ws.Range("Y5").Value = 1.5
$ws.Range("Z5").Value = 2.5
$ws.Range("AA5").Value = 1.91
$ws.Range("AB5").Value = 1.8
$ws.Range("AD5").Value = 10
$ws.Range("AE5").Value = 9.5
$ws.Range("AF5").Value = 21
$ws.Range("AG5").Value = 19
$ws.Range("AH5").Value = 34
$ws.Range("AI5").Value = 8.5
$ws.Range("AM5").Value = 8.5
$ws.Range("AN5").Value = 15
$ws.Range("AO5").Value = 12
$ws.Range("AP5").Value = 34
$ws.Range("AQ5").Value = 29
$ws.Range("AR5").Value = 41
$ws.Range("AS5").Value = 351
